$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CRMS-2314: add a "PAN" column to the vendor-list export template.
# Insert a brand-new column at N (right after "Phone 2 " / "Email"),
# pushing the existing "GST Number" ... "On/Off" columns one slot to
# the right.
$ws.Columns("N:N").Insert()

# Populate the two header rows for the newly inserted column. Row 2 is
# written before row 1 so the shared-string table records
# "{vendor:pan_no}" ahead of "PAN", matching how the columns were
# authored (merge-tag row, then display-label row).
$ws.Range("N2").Value = "{vendor:pan_no}"
$ws.Range("N1").Value = "PAN"

# Row 1 is the bold, centered header row used throughout the sheet -
# apply the same formatting to the new header cell.
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").HorizontalAlignment = -4108

# Leave the selection on the cell that used to be N1 (now shifted to
# O1), matching where the author's cursor ended up after inserting the
# column ahead of it.
[void]$ws.Range("O1").Select()
